$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 151 (the old blank separator row), shifting
# the blank row + summary rows down by one.
$ws.Rows("151:151").Insert()

# Fill in the new data row 151 with the same layout as rows 132-150
# (same formulas as the F/G columns above it, just resolved per-row).
$ws.Range("A151").Value = 2014
$ws.Range("B151").Value = 7
$ws.Range("C151").Value = 24
$ws.Range("D151").Value = 0.64583333333333337
$ws.Range("E151").Value = 0.75

$ws.Range("D151").NumberFormat = $ws.Range("D150").NumberFormat
$ws.Range("E151").NumberFormat = $ws.Range("E150").NumberFormat
$ws.Range("F151").NumberFormat = $ws.Range("F150").NumberFormat
$ws.Range("G151").NumberFormat = $ws.Range("G150").NumberFormat

$ws.Range("F132:F151").FormulaR1C1 = "=(RC[-1]-RC[-2])*24*60"
$ws.Range("G132:G151").FormulaR1C1 = "=RC[-1]/60"

# The totals block (now rows 153-155) needs its SUM() range extended to
# include the newly-inserted data row.
$ws.Range("F153").Formula = "=SUM(F2:F151)"

$ws.Range("F151").Select()
